# Add "2022-Q1" holdings sheet and register its totals on the "总计" sheet.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1 - Create the "2022-Q1" sheet.
# We copy "2021-Q4" because it already has the right 7-column layout/styles
# (基金代码/基金名称/基金规模/股票总仓位/仓位占比/持有市值(亿元)/仓位排名),
# and Excel's Copy places the new sheet immediately before the sheet we
# pass in - i.e. right before "总计", exactly where "2022-Q1" belongs.
# ---------------------------------------------------------------------------
$srcQ4 = $wb.Worksheets.Item("2021-Q4")
$total = $wb.Worksheets.Item("总计")
$srcQ4.Copy($total)

$q1 = $wb.Worksheets.Item("2021-Q4 (2)")
$q1.Name = "2022-Q1"
$q1 = $wb.Worksheets.Item("2022-Q1")

# The copied sheet has 3 data rows like "2021-Q4"; "2022-Q1" only has 2, so
# remove the extra (3rd) data row.
$q1.Rows.Item(4).Delete()
$q1 = $wb.Worksheets.Item("2022-Q1")

# Row 2 - rank 4 holding: fund 010764
$q1.Cells.Item(2,2).NumberFormat = "@"
$q1.Cells.Item(2,2).Value = "010764"
$q1.Cells.Item(2,2).ClearFormats()
$q1.Cells.Item(2,3).Value = "九泰锐升18个月封闭运作混合"
$q1.Cells.Item(2,4).NumberFormat = "@"
$q1.Cells.Item(2,4).Value = "3.15"
$q1.Cells.Item(2,4).ClearFormats()
$q1.Cells.Item(2,5).NumberFormat = "@"
$q1.Cells.Item(2,5).Value = "78.81"
$q1.Cells.Item(2,5).ClearFormats()
$q1.Cells.Item(2,6).NumberFormat = "@"
$q1.Cells.Item(2,6).Value = "2.93"
$q1.Cells.Item(2,6).ClearFormats()
$q1.Cells.Item(2,7).NumberFormat = "@"
$q1.Cells.Item(2,7).Value = "0.0923"
$q1.Cells.Item(2,7).ClearFormats()
$q1.Cells.Item(2,8).Value = 4

# Row 3 - rank 5 holding: fund 009531
$q1.Cells.Item(3,2).NumberFormat = "@"
$q1.Cells.Item(3,2).Value = "009531"
$q1.Cells.Item(3,2).ClearFormats()
$q1.Cells.Item(3,3).Value = "九泰锐和18个月定期开放混合"
$q1.Cells.Item(3,4).NumberFormat = "@"
$q1.Cells.Item(3,4).Value = "1.97"
$q1.Cells.Item(3,4).ClearFormats()
$q1.Cells.Item(3,5).NumberFormat = "@"
$q1.Cells.Item(3,5).Value = "70.75"
$q1.Cells.Item(3,5).ClearFormats()
$q1.Cells.Item(3,6).NumberFormat = "@"
$q1.Cells.Item(3,6).Value = "3.52"
$q1.Cells.Item(3,6).ClearFormats()
$q1.Cells.Item(3,7).NumberFormat = "@"
$q1.Cells.Item(3,7).Value = "0.0693"
$q1.Cells.Item(3,7).ClearFormats()
$q1.Cells.Item(3,8).Value = 5

# ---------------------------------------------------------------------------
# Step 2 - Add the 2022-Q1 summary row to the top of the "总计" sheet,
# pushing the existing quarters down and keeping the running index (column A)
# sequential.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

$total = $wb.Worksheets.Item("总计")
$total.Range("B2:D2").ClearFormats()
$total.Cells.Item(3,1).Copy()
$total.Cells.Item(2,1).PasteSpecial(-4122)
$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(2,2).Value = "2022-Q1"
$total.Cells.Item(2,3).Value = 2
$total.Cells.Item(2,4).Value = 0.16

# Renumber the remaining rows' index column (A) to stay sequential 0..3
$total.Cells.Item(3,1).Value = 1
$total.Cells.Item(4,1).Value = 2
$total.Cells.Item(5,1).Value = 3
